$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2886974440769166
$ws.Range("C2").Value = 0.04414062284175202
$ws.Range("E2").Value = 0.1519482541152186
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.6754222335128901
$ws.Range("H2").Value = 0.7796758169544518
$ws.Range("I2").Value = 0.6173599665478875
$ws.Range("K2").Value = 0.3061225245914443
$ws.Range("M2").Value = 0.2358896225223148
$ws.Range("N2").Value = 1.563560818350059
$ws.Range("B3").Value = 0.2572155069790654
$ws.Range("C3").Value = 0.0384478235663579
$ws.Range("E3").Value = 0.1409603973209244
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.6758928611237849
$ws.Range("H3").Value = 0.7842871477837008
$ws.Range("I3").Value = 0.6218221813652747
$ws.Range("K3").Value = 0.2706970287038075
$ws.Range("M3").Value = 0.2136452902220611
$ws.Range("N3").Value = 1.581963963130753
$ws.Range("B4").Value = 0.2379321969771127
$ws.Range("C4").Value = 0.03494855867897684
$ws.Range("E4").Value = 0.1343108316449886
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.6766622837825338
$ws.Range("H4").Value = 0.787488719127964
$ws.Range("I4").Value = 0.6249077154588534
$ws.Range("K4").Value = 0.2489808821832753
$ws.Range("M4").Value = 0.2000854270991823
$ws.Range("N4").Value = 1.593839392067274
$ws.Range("B5").Value = 0.230086112386374
$ws.Range("C5").Value = 0.03352158333005661
$ws.Range("E5").Value = 0.1316252766390349
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.6770963932726488
$ws.Range("H5").Value = 0.7888864510664391
$ws.Range("I5").Value = 0.6262519404952585
$ws.Range("K5").Value = 0.240140444407885
$ws.Range("M5").Value = 0.194584264002934
$ws.Range("N5").Value = 1.598823430780103
$ws.Range("B6").Value = 0.2287840100629239
$ws.Range("C6").Value = 0.03328457425298836
$ws.Range("E6").Value = 0.1311807978797219
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.6771757512371153
$ws.Range("H6").Value = 0.7891241636719499
$ws.Range("I6").Value = 0.6264803906822607
$ws.Range("K6").Value = 0.2386730502119292
$ws.Range("M6").Value = 0.1936722820280039
$ws.Range("N6").Value = 1.599659763832966
$ws.Range("B7").Value = 0.2378263329046888
$ws.Range("C7").Value = 0.03492931805888588
$ws.Range("E7").Value = 0.1342745156019802
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.6766676505118738
$ws.Range("H7").Value = 0.7875071926429058
$ws.Range("I7").Value = 0.6249254926347483
$ws.Range("K7").Value = 0.2488616199587739
$ws.Range("M7").Value = 0.2000111370697297
$ws.Range("N7").Value = 1.593906022774105
$ws.Range("B8").Value = 0.277832926568351
$ws.Range("C8").Value = 0.04217853062898769
$ws.Range("E8").Value = 0.1481393799301998
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.6754846342804939
$ws.Range("H8").Value = 0.7811889550514479
$ws.Range("I8").Value = 0.6188267311912874
$ws.Range("K8").Value = 0.2939006084600067
$ws.Range("M8").Value = 0.2281993099955386
$ws.Range("N8").Value = 1.569786657700623
$ws.Range("B9").Value = 0.3566485940698954
$ws.Range("C9").Value = 0.05636538753894627
$ws.Range("E9").Value = 0.1761086805017982
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.6769892239262276
$ws.Range("H9").Value = 0.7717376597618397
$ws.Range("I9").Value = 0.6096144253904896
$ws.Range("K9").Value = 0.3824971639359376
$ws.Range("M9").Value = 0.2842638681085887
$ws.Range("N9").Value = 1.527060864639921
$ws.Range("B10").Value = 0.4147709847960073
$ws.Range("C10").Value = 0.06677442317825921
$ws.Range("E10").Value = 0.1971501285520603
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.6804445448525911
$ws.Range("H10").Value = 0.7665878328729434
$ws.Range("I10").Value = 0.6045272546947373
$ws.Range("K10").Value = 0.4477574926536647
$ws.Range("M10").Value = 0.3259501322296501
$ws.Range("N10").Value = 1.498461026159147
$ws.Range("B11").Value = 0.4412587576086366
$ws.Range("C11").Value = 0.0715075314905107
$ws.Range("E11").Value = 0.2068330696744525
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.6825307778574938
$ws.Range("H11").Value = 0.7646352015366915
$ws.Range("I11").Value = 0.6025793953960346
$ws.Range("K11").Value = 0.4774832376856466
$ws.Range("M11").Value = 0.3450254989406005
$ws.Range("N11").Value = 1.486056970684935
$ws.Range("B12").Value = 0.4512956450120669
$ws.Range("C12").Value = 0.07329959041565814
$ws.Range("E12").Value = 0.2105159747467269
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.6833950540473239
$ws.Range("H12").Value = 0.7639519254378371
$ws.Range("I12").Value = 0.6018945795349495
$ws.Range("K12").Value = 0.4887450317470723
$ws.Range("M12").Value = 0.3522651543673945
$ws.Range("N12").Value = 1.481447149174096
$ws.Range("B13").Value = 0.4491337355565577
$ws.Range("C13").Value = 0.07291364994452465
$ws.Range("E13").Value = 0.2097220720091855
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.6832056090234317
$ws.Range("H13").Value = 0.7640965833218161
$ws.Range("I13").Value = 0.6020397168447715
$ws.Range("K13").Value = 0.4863193709475979
$ws.Range("M13").Value = 0.3507052415090683
$ws.Range("N13").Value = 1.482436070312067
$ws.Range("B14").Value = 0.4420843689842968
$ws.Range("C14").Value = 0.07165497046445068
$ws.Range("E14").Value = 0.2071357388942445
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.6826003922985535
$ws.Range("H14").Value = 0.7645778624976316
$ws.Range("I14").Value = 0.602521996380581
$ws.Range("K14").Value = 0.4784096473989905
$ws.Range("M14").Value = 0.3456207838746934
$ws.Range("N14").Value = 1.485675966837041
$ws.Range("B15").Value = 0.4377672714830396
$ws.Range("C15").Value = 0.0708839589290875
$ws.Range("E15").Value = 0.2055536493074683
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.6822393598831695
$ws.Range("H15").Value = 0.7648799729911531
$ws.Range("I15").Value = 0.6028242857085928
$ws.Range("K15").Value = 0.4735653965040285
$ws.Range("M15").Value = 0.3425085234528993
$ws.Range("N15").Value = 1.487671874142858
$ws.Range("B16").Value = 0.4130408785885606
$ws.Range("C16").Value = 0.06646506486497117
$ws.Range("E16").Value = 0.1965195789283243
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.6803185807841885
$ws.Range("H16").Value = 0.7667232953510421
$ws.Range("I16").Value = 0.6046619334045715
$ws.Range("K16").Value = 0.4458155997344022
$ws.Range("M16").Value = 0.3247057769464874
$ws.Range("N16").Value = 1.499283871403966
$ws.Range("B17").Value = 0.3978840146791924
$ws.Range("C17").Value = 0.06375371702860377
$ws.Range("E17").Value = 0.1910060705840948
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.6792722033668639
$ws.Range("H17").Value = 0.7679540535172578
$ws.Range("I17").Value = 0.6058831791640458
$ws.Range("K17").Value = 0.4288016895245619
$ws.Range("M17").Value = 0.3138131532900559
$ws.Range("N17").Value = 1.506562819680811
$ws.Range("B18").Value = 0.3891706863944933
$ws.Range("C18").Value = 0.06219403048230276
$ws.Range("E18").Value = 0.1878452891578632
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.6787187538325128
$ws.Range("H18").Value = 0.7686986631747885
$ws.Range("I18").Value = 0.6066200800052002
$ws.Range("K18").Value = 0.419019375694802
$ws.Range("M18").Value = 0.307558546704243
$ws.Range("N18").Value = 1.510806522494201
$ws.Range("B19").Value = 0.3862212850929723
$ws.Range("C19").Value = 0.06166591391493625
$ws.Range("E19").Value = 0.1867768899731104
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.6785396689013368
$ws.Range("H19").Value = 0.7689570785763209
$ws.Range("I19").Value = 0.6068754990302345
$ws.Range("K19").Value = 0.4157078842240765
$ws.Range("M19").Value = 0.3054426516602717
$ws.Range("N19").Value = 1.51225315992234
$ws.Range("B20").Value = 0.3994970252890084
$ws.Range("C20").Value = 0.06404236423259135
$ws.Range("E20").Value = 0.1915919106414776
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.6793785807592343
$ws.Range("H20").Value = 0.7678192374589088
$ws.Range("I20").Value = 0.6057496068177031
$ws.Range("K20").Value = 0.4306124745735929
$ws.Range("M20").Value = 0.3149715997647036
$ws.Range("N20").Value = 1.50578205818187
$ws.Range("B21").Value = 0.4441547636695873
$ws.Range("C21").Value = 0.07202468218457625
$ws.Range("E21").Value = 0.2078949666627707
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.6827761412249771
$ws.Range("H21").Value = 0.7644349749951687
$ws.Range("I21").Value = 0.6023789053591955
$ws.Range("K21").Value = 0.4807327829884969
$ws.Range("M21").Value = 0.347113771015124
$ws.Range("N21").Value = 1.484721959414914
$ws.Range("B22").Value = 0.4733791729411507
$ws.Range("C22").Value = 0.07724004608712676
$ws.Range("E22").Value = 0.2186443960874129
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.6854296547511893
$ws.Range("H22").Value = 0.7625504341119296
$ws.Range("I22").Value = 0.6004837496867594
$ws.Range("K22").Value = 0.513520082653713
$ws.Range("M22").Value = 0.3682151795337703
$ws.Range("N22").Value = 1.471467041918403
$ws.Range("B23").Value = 0.4577781790568736
$ws.Range("C23").Value = 0.07445664163880394
$ws.Range("E23").Value = 0.2128985117975333
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.6839737060122815
$ws.Range("H23").Value = 0.7635262856418024
$ws.Range("I23").Value = 0.6014670271554934
$ws.Range("K23").Value = 0.4960181451171479
$ws.Range("M23").Value = 0.3569442639856248
$ws.Range("N23").Value = 1.478494806494691
$ws.Range("B24").Value = 0.3987677814420181
$ws.Range("C24").Value = 0.06391186962643758
$ws.Range("E24").Value = 0.1913270243118532
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.6793303376481816
$ws.Range("H24").Value = 0.7678800724518737
$ws.Range("I24").Value = 0.6058098865118069
$ws.Range("K24").Value = 0.4297938210986842
$ws.Range("M24").Value = 0.3144478421167349
$ws.Range("N24").Value = 1.506134856869513
$ws.Range("B25").Value = 0.3352885085284072
$ws.Range("C25").Value = 0.05253018924254604
$ws.Range("E25").Value = 0.1684569799475355
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.6761708571439584
$ws.Range("H25").Value = 0.7739796178343283
$ws.Range("I25").Value = 0.6118118320530037
$ws.Range("K25").Value = 0.3584999179147133
$ws.Range("M25").Value = 0.2690110002537125
$ws.Range("N25").Value = 1.538129492881605